# Market name on sheet
# - Rename the sheet to "VARIOS"
# - Update column A (price) values for a set of rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "VARIOS"

$ws.Cells.Item(2, 1).Value = 15799
$ws.Cells.Item(3, 1).Value = 12199
$ws.Cells.Item(4, 1).Value = 6299
$ws.Cells.Item(5, 1).Value = 8999
$ws.Cells.Item(6, 1).Value = 17599
$ws.Cells.Item(7, 1).Value = 10799
$ws.Cells.Item(8, 1).Value = 14799
$ws.Cells.Item(9, 1).Value = 15699
$ws.Cells.Item(10, 1).Value = 16299
$ws.Cells.Item(11, 1).Value = 64999
$ws.Cells.Item(12, 1).Value = 20799
$ws.Cells.Item(13, 1).Value = 24299
$ws.Cells.Item(14, 1).Value = 22999
$ws.Cells.Item(15, 1).Value = 29899
$ws.Cells.Item(16, 1).Value = 6499
$ws.Cells.Item(17, 1).Value = 4999
$ws.Cells.Item(18, 1).Value = 18499
$ws.Cells.Item(19, 1).Value = 13199
$ws.Cells.Item(20, 1).Value = 12699
$ws.Cells.Item(21, 1).Value = 34199
$ws.Cells.Item(22, 1).Value = 14299
$ws.Cells.Item(23, 1).Value = 11999
$ws.Cells.Item(24, 1).Value = 8299
$ws.Cells.Item(25, 1).Value = 10299
$ws.Cells.Item(26, 1).Value = 14099
$ws.Cells.Item(27, 1).Value = 9099
$ws.Cells.Item(28, 1).Value = 29399
$ws.Cells.Item(29, 1).Value = 19799
$ws.Cells.Item(30, 1).Value = 8599
$ws.Cells.Item(31, 1).Value = 42099
$ws.Cells.Item(32, 1).Value = 24899
$ws.Cells.Item(33, 1).Value = 15899
$ws.Cells.Item(34, 1).Value = 25599
$ws.Cells.Item(35, 1).Value = 20499
$ws.Cells.Item(37, 1).Value = 16199
